$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 6
$ws_ALC.Range("H6").Value = 1000
$ws_ALC.Range("J6").Value = 1000
$ws_ALC.Range("L6").Value = 3000
$ws_ALC.Range("N6").Value = -3224

# ALC row 31
$ws_ALC.Range("H31").Value = 2935.3333
$ws_ALC.Range("I31").Value = 903
$ws_ALC.Range("K31").Value = 2709
$ws_ALC.Range("M31").Value = -2479

# ALC row 62
$ws_ALC.Range("H62").Value = 2909.5
$ws_ALC.Range("I62").Value = 2428.375
$ws_ALC.Range("J62").Value = 4064.2
$ws_ALC.Range("K62").Value = 2428.375
$ws_ALC.Range("L62").Value = 4064.2
$ws_ALC.Range("M62").Value = -1804.375
$ws_ALC.Range("N62").Value = -5312.2

# ALC row 65
$ws_ALC.Range("H65").Value = 2909.5
$ws_ALC.Range("I65").Value = 2428.375
$ws_ALC.Range("J65").Value = 4064.2
$ws_ALC.Range("K65").Value = 12141.875
$ws_ALC.Range("L65").Value = 20321
$ws_ALC.Range("M65").Value = -9021.875
$ws_ALC.Range("N65").Value = -26561

# ALC row 82
$ws_ALC.Range("H82").Value = 538.5
$ws_ALC.Range("I82").Value = 538.5
$ws_ALC.Range("K82").Value = 1615.5
$ws_ALC.Range("M82").Value = -1209.5

# ALC row 85
$ws_ALC.Range("H85").Value = 538.5
$ws_ALC.Range("I85").Value = 538.5
$ws_ALC.Range("K85").Value = 1615.5
$ws_ALC.Range("M85").Value = -211.5

# ALC row 86
$ws_ALC.Range("H86").Value = 13393.25
$ws_ALC.Range("I86").Value = 1180
$ws_ALC.Range("K86").Value = 1180
$ws_ALC.Range("M86").Value = -57

# ALC row 88
$ws_ALC.Range("H88").Value = 1087.5
$ws_ALC.Range("I88").Value = 675
$ws_ALC.Range("J88").Value = 1500
$ws_ALC.Range("K88").Value = 675
$ws_ALC.Range("L88").Value = 1500
$ws_ALC.Range("M88").Value = -269
$ws_ALC.Range("N88").Value = -2312

# ALC row 89
$ws_ALC.Range("H89").Value = 13393.25
$ws_ALC.Range("I89").Value = 1180
$ws_ALC.Range("K89").Value = 5900
$ws_ALC.Range("M89").Value = -284

# ALC row 91
$ws_ALC.Range("H91").Value = 1087.5
$ws_ALC.Range("I91").Value = 675
$ws_ALC.Range("J91").Value = 1500
$ws_ALC.Range("K91").Value = 675
$ws_ALC.Range("L91").Value = 1500
$ws_ALC.Range("M91").Value = 729
$ws_ALC.Range("N91").Value = -4308

# ALC row 93
$ws_ALC.Range("H93").Value = 28600
$ws_ALC.Range("J93").Value = 28600
$ws_ALC.Range("L93").Value = 28600
$ws_ALC.Range("N93").Value = -33592

# ALC row 113
$ws_ALC.Range("H113").Value = 4622.4443
$ws_ALC.Range("J113").Value = 4600.3335
$ws_ALC.Range("L113").Value = 4600.3335
$ws_ALC.Range("N113").Value = -11108.3335

# ALC row 129
$ws_ALC.Range("H129").Value = 347362.28
$ws_ALC.Range("J129").Value = 359750.22
$ws_ALC.Range("L129").Value = 1079250.66
$ws_ALC.Range("N129").Value = -1089250.66

# ALC row 132
$ws_ALC.Range("H132").Value = 3267.7693
$ws_ALC.Range("I132").Value = 3267.7693
$ws_ALC.Range("K132").Value = 9803.3079
$ws_ALC.Range("M132").Value = -7273.3079

# ALC row 141
$ws_ALC.Range("H141").Value = 1766.8438
$ws_ALC.Range("I141").Value = 1221.5
$ws_ALC.Range("J141").Value = 2675.75
$ws_ALC.Range("K141").Value = 3664.5
$ws_ALC.Range("L141").Value = 8027.25
$ws_ALC.Range("M141").Value = 1515.5
$ws_ALC.Range("N141").Value = -18387.25

# ARM row 59
$ws_ARM.Range("H59").Value = 22000
$ws_ARM.Range("J59").Value = 22000
$ws_ARM.Range("L59").Value = 22000
$ws_ARM.Range("N59").Value = -23608

# ARM row 61
$ws_ARM.Range("H61").Value = 1261.9678
$ws_ARM.Range("I61").Value = 1319.7037
$ws_ARM.Range("J61").Value = 872.25
$ws_ARM.Range("K61").Value = 1319.7037
$ws_ARM.Range("L61").Value = 872.25
$ws_ARM.Range("M61").Value = -1107.7037
$ws_ARM.Range("N61").Value = -1296.25

# ARM row 122
$ws_ARM.Range("H122").Value = 2111.3845
$ws_ARM.Range("I122").Value = 1588
$ws_ARM.Range("K122").Value = 4764
$ws_ARM.Range("M122").Value = -2314

# ARM row 136
$ws_ARM.Range("H136").Value = 1261.9678
$ws_ARM.Range("I136").Value = 1319.7037
$ws_ARM.Range("J136").Value = 872.25
$ws_ARM.Range("K136").Value = 3959.1111
$ws_ARM.Range("L136").Value = 2616.75
$ws_ARM.Range("M136").Value = -1409.1111
$ws_ARM.Range("N136").Value = -7716.75

# BSM row 86
$ws_BSM.Range("H86").Value = 1371.579
$ws_BSM.Range("I86").Value = 1261.6
$ws_BSM.Range("J86").Value = 1583.0769
$ws_BSM.Range("K86").Value = 1261.6
$ws_BSM.Range("L86").Value = 1583.0769
$ws_BSM.Range("M86").Value = -138.5999999999999
$ws_BSM.Range("N86").Value = -3829.0769

# BSM row 89
$ws_BSM.Range("H89").Value = 1371.579
$ws_BSM.Range("I89").Value = 1261.6
$ws_BSM.Range("J89").Value = 1583.0769
$ws_BSM.Range("K89").Value = 6308
$ws_BSM.Range("L89").Value = 7915.3845
$ws_BSM.Range("M89").Value = -692
$ws_BSM.Range("N89").Value = -19147.3845

# CRP row 31
$ws_CRP.Range("H31").Value = 8834.851000000001
$ws_CRP.Range("I31").Value = 10118.229
$ws_CRP.Range("K31").Value = 10118.229
$ws_CRP.Range("M31").Value = -9823.228999999999

# CRP row 34
$ws_CRP.Range("H34").Value = 8834.851000000001
$ws_CRP.Range("I34").Value = 10118.229
$ws_CRP.Range("K34").Value = 10118.229
$ws_CRP.Range("M34").Value = -9916.228999999999

# CUL row 86
$ws_CUL.Range("H86").Value = 71429020
$ws_CUL.Range("I86").Value = 573
$ws_CUL.Range("K86").Value = 1719
$ws_CUL.Range("M86").Value = -533

# CUL row 89
$ws_CUL.Range("H89").Value = 71429020
$ws_CUL.Range("I89").Value = 573
$ws_CUL.Range("K89").Value = 5157
$ws_CUL.Range("M89").Value = 771

# CUL row 122
$ws_CUL.Range("H122").Value = 765.5
$ws_CUL.Range("J122").Value = 858.75
$ws_CUL.Range("L122").Value = 7728.75
$ws_CUL.Range("N122").Value = -12628.75

# CUL row 131
$ws_CUL.Range("H131").Value = 121333.02
$ws_CUL.Range("J131").Value = 132459.62
$ws_CUL.Range("L131").Value = 397378.86
$ws_CUL.Range("N131").Value = -407458.86

# GSM row 80
$ws_GSM.Range("H80").Value = 3512.2083
$ws_GSM.Range("I80").Value = 2817.5454
$ws_GSM.Range("K80").Value = 2817.5454
$ws_GSM.Range("M80").Value = -1819.5454

# GSM row 83
$ws_GSM.Range("H83").Value = 3512.2083
$ws_GSM.Range("I83").Value = 2817.5454
$ws_GSM.Range("K83").Value = 14087.727
$ws_GSM.Range("M83").Value = -9095.726999999999

# GSM row 92
$ws_GSM.Range("H92").Value = 16666.666
$ws_GSM.Range("I92").Value = 0
$ws_GSM.Range("J92").Value = 16666.666
$ws_GSM.Range("K92").Value = 0
$ws_GSM.Range("L92").Value = 16666.666
$ws_GSM.Range("M92").ClearContents()
$ws_GSM.Range("N92").Value = -20410.666

# GSM row 126
$ws_GSM.Range("H126").Value = 3741.1628
$ws_GSM.Range("I126").Value = 2980.2
$ws_GSM.Range("J126").Value = 5497.231
$ws_GSM.Range("K126").Value = 8940.599999999999
$ws_GSM.Range("L126").Value = 16491.693
$ws_GSM.Range("M126").Value = -6470.599999999999
$ws_GSM.Range("N126").Value = -21431.693

# LTW row 55
$ws_LTW.Range("H55").Value = 295.625
$ws_LTW.Range("I55").Value = 391.25
$ws_LTW.Range("J55").Value = 200
$ws_LTW.Range("K55").Value = 391.25
$ws_LTW.Range("L55").Value = 200
$ws_LTW.Range("M55").Value = -218.25
$ws_LTW.Range("N55").Value = -546

# LTW row 61
$ws_LTW.Range("H61").Value = 7500
$ws_LTW.Range("I61").Value = 3833.3333
$ws_LTW.Range("J61").Value = 10250
$ws_LTW.Range("K61").Value = 3833.3333
$ws_LTW.Range("L61").Value = 10250
$ws_LTW.Range("M61").Value = -3631.3333
$ws_LTW.Range("N61").Value = -10654

# LTW row 113
$ws_LTW.Range("H113").Value = 7500
$ws_LTW.Range("I113").Value = 3833.3333
$ws_LTW.Range("J113").Value = 10250
$ws_LTW.Range("K113").Value = 3833.3333
$ws_LTW.Range("L113").Value = 10250
$ws_LTW.Range("M113").Value = -1663.3333
$ws_LTW.Range("N113").Value = -14590

# WVR row 107
$ws_WVR.Range("H107").Value = 3031051.5
$ws_WVR.Range("J107").Value = 6494091.5
$ws_WVR.Range("L107").Value = 19482274.5
$ws_WVR.Range("N107").Value = -19486114.5
